$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pop_names")
$ws2 = $wb.Worksheets.Item("Colors")

# --- Colors sheet: reorder rows 1-13 and add Order column C ---
$ws2.Cells.Item(1,1).Value = "African Dogs"
$ws2.Cells.Item(1,2).Value = "#a6cee3"
$ws2.Cells.Item(1,3).Value = 1
$ws2.Cells.Item(2,1).Value = "European Dogs"
$ws2.Cells.Item(2,2).Value = "#33a02c"
$ws2.Cells.Item(2,3).Value = 2
$ws2.Cells.Item(3,1).Value = "Northern Dogs"
$ws2.Cells.Item(3,2).Value = "#fb9a99"
$ws2.Cells.Item(3,3).Value = 3
$ws2.Cells.Item(4,1).Value = "American Dogs"
$ws2.Cells.Item(4,2).Value = "#cab2d6"
$ws2.Cells.Item(4,3).Value = 4
$ws2.Cells.Item(5,1).Value = "Asian Dogs"
$ws2.Cells.Item(5,2).Value = "#fdbf6f"
$ws2.Cells.Item(5,3).Value = 5
$ws2.Cells.Item(6,1).Value = "East Asian Dogs"
$ws2.Cells.Item(6,2).Value = "#e31a1c"
$ws2.Cells.Item(6,3).Value = 6
$ws2.Cells.Item(7,1).Value = "Pre-Colombian Dogs"
$ws2.Cells.Item(7,2).Value = "#6a3d9a"
$ws2.Cells.Item(7,3).Value = 7
$ws2.Cells.Item(8,1).Value = "CTVT"
$ws2.Cells.Item(8,2).Value = "#b2df8a"
$ws2.Cells.Item(8,3).Value = 8
$ws2.Cells.Item(9,1).Value = "Dingo"
$ws2.Cells.Item(9,2).Value = "#ff7f00"
$ws2.Cells.Item(9,3).Value = 9
$ws2.Cells.Item(10,1).Value = "Coyotes"
$ws2.Cells.Item(10,2).Value = "#1f78b4"
$ws2.Cells.Item(10,3).Value = 10
$ws2.Cells.Item(11,1).Value = "American Wolf"
$ws2.Cells.Item(11,2).Value = "#b15928"
$ws2.Cells.Item(11,3).Value = 11
$ws2.Cells.Item(12,1).Value = "Eurasian Wolf"
$ws2.Cells.Item(12,2).Value = "#003c30"
$ws2.Cells.Item(12,3).Value = 12
$ws2.Cells.Item(13,1).Value = "Outgroup"
$ws2.Cells.Item(13,2).Value = "#ffff99"
$ws2.Cells.Item(13,3).Value = 13

# --- pop_names sheet: reorder rows 2-29, add Order header + column E formula ---
$ws1.Cells.Item(1,5).Value = "Order"
$ws1.Cells.Item(2,1).Value = "BAS"
$ws1.Cells.Item(2,2).Value = "Basenji"
$ws1.Cells.Item(2,3).Value = "African Dogs"
$ws1.Cells.Item(2,4).Formula = "=VLOOKUP(C2,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(2,5).Formula = "=VLOOKUP(D2,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(3,1).Value = "DNA"
$ws1.Cells.Item(3,2).Value = "Namibian Village Dog"
$ws1.Cells.Item(3,3).Value = "African Dogs"
$ws1.Cells.Item(3,4).Formula = "=VLOOKUP(C3,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(3,5).Formula = "=VLOOKUP(D3,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(4,1).Value = "DAE"
$ws1.Cells.Item(4,2).Value = "Ancient European"
$ws1.Cells.Item(4,3).Value = "European Dogs"
$ws1.Cells.Item(4,4).Formula = "=VLOOKUP(C4,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(4,5).Formula = "=VLOOKUP(D4,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(5,1).Value = "DEU"
$ws1.Cells.Item(5,2).Value = "European Village Dog"
$ws1.Cells.Item(5,3).Value = "European Dogs"
$ws1.Cells.Item(5,4).Formula = "=VLOOKUP(C5,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(5,5).Formula = "=VLOOKUP(D5,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(6,1).Value = "DGS"
$ws1.Cells.Item(6,2).Value = "German Shepard"
$ws1.Cells.Item(6,3).Value = "European Dogs"
$ws1.Cells.Item(6,4).Formula = "=VLOOKUP(C6,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(6,5).Formula = "=VLOOKUP(D6,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(7,1).Value = "DLB"
$ws1.Cells.Item(7,2).Value = "Lebanese Village Dog"
$ws1.Cells.Item(7,3).Value = "European Dogs"
$ws1.Cells.Item(7,4).Formula = "=VLOOKUP(C7,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(7,5).Formula = "=VLOOKUP(D7,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(8,1).Value = "DAL"
$ws1.Cells.Item(8,2).Value = "Alaskan Husky"
$ws1.Cells.Item(8,3).Value = "Northern Dogs"
$ws1.Cells.Item(8,4).Formula = "=VLOOKUP(C8,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(8,5).Formula = "=VLOOKUP(D8,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(9,1).Value = "DGL"
$ws1.Cells.Item(9,2).Value = "Greenland Sledge Dog"
$ws1.Cells.Item(9,3).Value = "Northern Dogs"
$ws1.Cells.Item(9,4).Formula = "=VLOOKUP(C9,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(9,5).Formula = "=VLOOKUP(D9,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(10,1).Value = "DHU"
$ws1.Cells.Item(10,2).Value = "Husky"
$ws1.Cells.Item(10,3).Value = "Northern Dogs"
$ws1.Cells.Item(10,4).Formula = "=VLOOKUP(C10,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(10,5).Formula = "=VLOOKUP(D10,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(11,1).Value = "DMA"
$ws1.Cells.Item(11,2).Value = "Malamute"
$ws1.Cells.Item(11,3).Value = "Northern Dogs"
$ws1.Cells.Item(11,4).Formula = "=VLOOKUP(C11,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(11,5).Formula = "=VLOOKUP(D11,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(12,1).Value = "DSL"
$ws1.Cells.Item(12,2).Value = "Siberian Laika"
$ws1.Cells.Item(12,3).Value = "Northern Dogs"
$ws1.Cells.Item(12,4).Formula = "=VLOOKUP(C12,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(12,5).Formula = "=VLOOKUP(D12,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(13,1).Value = "DME"
$ws1.Cells.Item(13,2).Value = "Mexican Hairless Dog"
$ws1.Cells.Item(13,3).Value = "American Dogs"
$ws1.Cells.Item(13,4).Formula = "=VLOOKUP(C13,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(13,5).Formula = "=VLOOKUP(D13,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(14,1).Value = "DPU"
$ws1.Cells.Item(14,2).Value = "Peruvian Hairless Dog"
$ws1.Cells.Item(14,3).Value = "American Dogs"
$ws1.Cells.Item(14,4).Formula = "=VLOOKUP(C14,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(14,5).Formula = "=VLOOKUP(D14,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(15,1).Value = "DID"
$ws1.Cells.Item(15,2).Value = "Indian Village Dog"
$ws1.Cells.Item(15,3).Value = "Asian Dogs"
$ws1.Cells.Item(15,4).Formula = "=VLOOKUP(C15,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(15,5).Formula = "=VLOOKUP(D15,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(16,1).Value = "DQA"
$ws1.Cells.Item(16,2).Value = "Qatari Village Dogs"
$ws1.Cells.Item(16,3).Value = "Asian Dogs"
$ws1.Cells.Item(16,4).Formula = "=VLOOKUP(C16,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(16,5).Formula = "=VLOOKUP(D16,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(17,1).Value = "DCH"
$ws1.Cells.Item(17,2).Value = "Chinese Village Dog"
$ws1.Cells.Item(17,3).Value = "East Asian Dogs"
$ws1.Cells.Item(17,4).Formula = "=VLOOKUP(C17,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(17,5).Formula = "=VLOOKUP(D17,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(18,1).Value = "DTI"
$ws1.Cells.Item(18,2).Value = "Tibetan Village Dog"
$ws1.Cells.Item(18,3).Value = "East Asian Dogs"
$ws1.Cells.Item(18,4).Formula = "=VLOOKUP(C18,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(18,5).Formula = "=VLOOKUP(D18,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(19,1).Value = "DTM"
$ws1.Cells.Item(19,2).Value = "Tibetan Mastiff"
$ws1.Cells.Item(19,3).Value = "East Asian Dogs"
$ws1.Cells.Item(19,4).Formula = "=VLOOKUP(C19,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(19,5).Formula = "=VLOOKUP(D19,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(20,1).Value = "DVN"
$ws1.Cells.Item(20,2).Value = "Vietnamese Village Dog"
$ws1.Cells.Item(20,3).Value = "East Asian Dogs"
$ws1.Cells.Item(20,4).Formula = "=VLOOKUP(C20,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(20,5).Formula = "=VLOOKUP(D20,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(21,1).Value = "DPC"
$ws1.Cells.Item(21,2).Value = "Pre-Colombian Dogs"
$ws1.Cells.Item(21,3).Value = "Pre-Colombian Dogs"
$ws1.Cells.Item(21,4).Formula = "=VLOOKUP(C21,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(21,5).Formula = "=VLOOKUP(D21,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(22,1).Value = "CTVT"
$ws1.Cells.Item(22,2).Value = "CTVT"
$ws1.Cells.Item(22,3).Value = "CTVT"
$ws1.Cells.Item(22,4).Formula = "=VLOOKUP(C22,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(22,5).Formula = "=VLOOKUP(D22,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(23,1).Value = "DIN"
$ws1.Cells.Item(23,2).Value = "Dingo"
$ws1.Cells.Item(23,3).Value = "Dingo"
$ws1.Cells.Item(23,4).Formula = "=VLOOKUP(C23,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(23,5).Formula = "=VLOOKUP(D23,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(24,1).Value = "COY"
$ws1.Cells.Item(24,2).Value = "Coyote"
$ws1.Cells.Item(24,3).Value = "Coyotes"
$ws1.Cells.Item(24,4).Formula = "=VLOOKUP(C24,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(24,5).Formula = "=VLOOKUP(D24,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(25,1).Value = "WAM"
$ws1.Cells.Item(25,2).Value = "American Wolf"
$ws1.Cells.Item(25,3).Value = "American Wolf"
$ws1.Cells.Item(25,4).Formula = "=VLOOKUP(C25,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(25,5).Formula = "=VLOOKUP(D25,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(26,1).Value = "WAS"
$ws1.Cells.Item(26,2).Value = "Asian Wolf"
$ws1.Cells.Item(26,3).Value = "Eurasian Wolf"
$ws1.Cells.Item(26,4).Formula = "=VLOOKUP(C26,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(26,5).Formula = "=VLOOKUP(D26,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(27,1).Value = "WEU"
$ws1.Cells.Item(27,2).Value = "European Wolf"
$ws1.Cells.Item(27,3).Value = "Eurasian Wolf"
$ws1.Cells.Item(27,4).Formula = "=VLOOKUP(C27,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(27,5).Formula = "=VLOOKUP(D27,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(28,1).Value = "WME"
$ws1.Cells.Item(28,2).Value = "Middle-east Wolf"
$ws1.Cells.Item(28,3).Value = "Eurasian Wolf"
$ws1.Cells.Item(28,4).Formula = "=VLOOKUP(C28,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(28,5).Formula = "=VLOOKUP(D28,Colors!B:C,2,FALSE)"
$ws1.Cells.Item(29,1).Value = "OUT"
$ws1.Cells.Item(29,2).Value = "Andean Fox"
$ws1.Cells.Item(29,3).Value = "Outgroup"
$ws1.Cells.Item(29,4).Formula = "=VLOOKUP(C29,Colors!A:B,2,FALSE)"
$ws1.Cells.Item(29,5).Formula = "=VLOOKUP(D29,Colors!B:C,2,FALSE)"

# --- AutoFilter + FilterDatabase name ---
$ws1.Range("A1:E29").AutoFilter()
$n = $ws1.Names.Add("_xlnm._FilterDatabase", "=pop_names!`$A`$1:`$E`$29")
$n.Visible = $false

# --- selection cells to match authored view state ---
$ws1.Range("G3").Select()
$ws2.Range("E10").Select()
